# Swap the deck's two theme colour schemes: the slide master/presentation
# theme (theme1.xml, currently "Integral"/"Red Violet") becomes the plain
# "Office" colour scheme that previously only lived in the notes-master
# theme (theme2.xml). Font scheme / format scheme are already identical
# between the two theme parts, so only the colour scheme needs editing.
#
# Per this host's own guidance, theme colours are edited through
# ThemeColorScheme.Item(i).RGB (VBA-style BGR long), not by touching the
# underlying theme XML part directly.

$p = $ppt.ActivePresentation
$t = $p.SlideMaster.Theme
$tcs = $t.ThemeColorScheme

# Office theme colour scheme, in slot order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# (values are the RGB long -- 0xBBGGRR -- that COM's ColorFormat.RGB uses)
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
